# surfspots.xlsx - "added new spots and debugged"
# Adds 5 new surf-spot rows (14-18) to the spot_database sheet.
# Columns: A=nomSpot, B=villeSpot, C=nomSurfForecast, D=paysSpot
# Cells are written in the same per-row order the spreadsheet was
# originally authored in (forecast name, then spot name, then city,
# then country) so the shared-string table comes out identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - A Corunna (Spain) spot
$ws.Cells.Item(14, 3).Value = "Playade-Bastiagueiros"
$ws.Cells.Item(14, 1).Value = "Platja de Bastiagueiro"
$ws.Cells.Item(14, 2).Value = "La Corogne"
$ws.Cells.Item(14, 4).Value = "Espagne"

# Row 15 - Nice
$ws.Cells.Item(15, 3).Value = "La-Marina-1"
$ws.Cells.Item(15, 2).Value = "Nice"
$ws.Cells.Item(15, 1).Value = "La Marina"
$ws.Cells.Item(15, 4).Value = "France"

# Row 16 - Biarritz
$ws.Cells.Item(16, 1).Value = "Biarritz Grande-Plage"
$ws.Cells.Item(16, 3).Value = "Grande-Plage"
$ws.Cells.Item(16, 2).Value = "Biarritz"
$ws.Cells.Item(16, 4).Value = "France"

# Row 17 - Mimizan
$ws.Cells.Item(17, 1).Value = "Mimizan"
$ws.Cells.Item(17, 2).Value = "Mimizan"
$ws.Cells.Item(17, 3).Value = "Mimizan"
$ws.Cells.Item(17, 4).Value = "France"

# Row 18 - Pornichet
$ws.Cells.Item(18, 1).Value = "Pornichet"
$ws.Cells.Item(18, 2).Value = "Pornichet"
$ws.Cells.Item(18, 3).Value = "Pornichet"
$ws.Cells.Item(18, 4).Value = "France"

Write-Output "Added 5 new surf spots (rows 14-18)"
